$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Races")

$ws.Range("D9").Value = "Verstappen"
$ws.Range("D10").Value = "Verstappen"
$ws.Range("D11").Value = "Redbull"
$ws.Range("D12").Value = "Leclerc"
$ws.Range("D13").Value = "Zhou"
